$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (Jul 22 2024 GitHub Actions refresh)

$ws.Range("D2").Value = "'67.408.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.02%  '

$ws.Range("D3").Value = "'3.435.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.73%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = "'589.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.10%  '

$ws.Range("D6").Value = "'178.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.51%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = "'0.602"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.86%  '

$ws.Range("D9").Value = "'3.428.94"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.90%  '

$ws.Range("D10").Value = "'0.137"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.36%  '

$ws.Range("D11").Value = "'6.95"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.87%  '

$ws.Range("D12").Value = "'0.425"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.01%  '

$ws.Range("D13").Value = "'4.032.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.76%  '

$ws.Range("D14").Value = "'31.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.61%  '

$ws.Range("E15").Value = '  -1.23%  '

$ws.Range("D16").Value = "'67.382.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.98%  '

$ws.Range("E17").Value = '  -3.90%  '

$ws.Range("D18").Value = "'3.435.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.70%  '

$ws.Range("D19").Value = "'6.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.05%  '

$ws.Range("D20").Value = "'13.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -6.96%  '

$ws.Range("D21").Value = "'384.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.01%  '

$ws.Range("D22").Value = "'7.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.60%  '

$ws.Range("D23").Value = "'5.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.55%  '

$ws.Range("D24").Value = "'0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.18%  '

$ws.Range("B25").Value = 'Polygon'
$ws.Range("C25").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D25").Value = "'0.529"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.56%  '

$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").Value = "'70.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.10%  '

$ws.Range("D27").Value = "'0.0000118"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.48%  '

$ws.Range("E28").Value = '  -5.22%  '

$ws.Range("E29").Value = '  -2.89%  '

$ws.Range("E30").Value = '  +0.34%  '

$ws.Range("E31").Value = '  -4.82%  '

$ws.Range("D32").Value = "'2.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.66%  '

$ws.Range("E33").Value = '  -7.61%  '

$ws.Range("D34").Value = "'23.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.34%  '

$ws.Range("E35").Value = '  -0.07%  '

$ws.Range("D36").Value = "'7.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.91%  '

$ws.Range("D37").Value = "'1.54"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -8.18%  '

$ws.Range("D38").Value = "'160.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.89%  '

$ws.Range("D39").Value = "'0.877"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.73%  '

$ws.Range("D40").Value = "'1.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.89%  '

$ws.Range("E41").Value = '  -3.94%  '

$ws.Range("D42").Value = "'6.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.64%  '

$ws.Range("D43").Value = "'4.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.03%  '

$ws.Range("D44").Value = "'25.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.23%  '

$ws.Range("D45").Value = "'0.0708"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.97%  '

$ws.Range("D46").Value = "'25.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.77%  '

$ws.Range("D47").Value = "'2.691.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.73%  '

$ws.Range("D48").Value = "'41.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.44%  '

$ws.Range("D49").Value = "'0.0294"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.22%  '

$ws.Range("D50").Value = "'323.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -8.32%  '

$ws.Range("D51").Value = "'1.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.26%  '
